$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $value) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2" "307.55"
Set-TextValue "E2" "-2.53%"

Set-TextValue "D3" "37.70"
Set-TextValue "E3" "-4.42%"

Set-TextValue "D4" "5.098"
Set-TextValue "E4" "-0.34%"

Set-TextValue "D5" "0.07870"
Set-TextValue "E5" "-3.91%"

Set-TextValue "D6" "1.974"
Set-TextValue "E6" "-3.23%"

Set-TextValue "D7" "4.341"
Set-TextValue "E7" "1.54%"

Set-TextValue "D8" "8.246"
Set-TextValue "E8" "-0.15%"

Set-TextValue "E9" "-5.76%"

Set-TextValue "D10" "0.9306"
Set-TextValue "E10" "-0.23%"

Set-TextValue "D11" "0.1299"
Set-TextValue "E11" "-7.90%"

Set-TextValue "D12" "0.1912"
Set-TextValue "E12" "-4.00%"

Set-TextValue "D13" "0.08869"
Set-TextValue "E13" "-2.78%"

Set-TextValue "D14" "0.03428"

Set-TextValue "D15" "0.09745"
Set-TextValue "E15" "-0.63%"

Set-TextValue "D16" "0.001396"
Set-TextValue "E16" "-0.12%"

Set-TextValue "D17" "0.005894"
Set-TextValue "E17" "-6.77%"

Set-TextValue "E18" "1,775.85%"

Set-TextValue "D19" "3.564"
Set-TextValue "E19" "-2.51%"

Set-TextValue "D20" "0.3433"
Set-TextValue "E20" "-0.76%"

Set-TextValue "E21" "-1.21%"

Set-TextValue "E22" "1.99%"

Set-TextValue "D23" "0.2487"
Set-TextValue "E23" "1.50%"

Set-TextValue "D24" "0.04315"
Set-TextValue "E24" "-0.34%"

Set-TextValue "D25" "0.001217"
Set-TextValue "E25" "-0.63%"

Set-TextValue "D26" "0.004607"
Set-TextValue "E26" "-3.56%"

Set-TextValue "E27" "176.18%"

Set-TextValue "D39" "0.02328"
Set-TextValue "E39" "3.67%"

Set-TextValue "D40" "0.05029"
Set-TextValue "E40" "-4.33%"

Set-TextValue "D41" "0.007501"
Set-TextValue "E41" "-0.21%"

Set-TextValue "D42" "0.009765"
Set-TextValue "E42" "-0.22%"

Set-TextValue "D43" "0.1354"
Set-TextValue "E43" "-1.58%"

Set-TextValue "D44" "0.002090"
Set-TextValue "E44" "-2.83%"

Set-TextValue "E45" "-15.41%"

Set-TextValue "D46" "0.00006557"
Set-TextValue "E46" "1.62%"

Set-TextValue "E47" "-0.07%"

Set-TextValue "E48" "8.24%"

Set-TextValue "E50" "-0.07%"

Set-TextValue "E51" "-0.07%"
